$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column for the "ssim_dual" experiment results,
# shifting the existing ssim_nlm / ssim_gnlm / ssim_bm3d columns right.
$ws.Columns("A").Insert()

# Header
$ws.Range("A1").Value = "ssim_dual"
# Match the header formatting used by the other header cells
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# New per-image ssim_dual values (rows 2-51)
$values = @(0.58989563, 0.67001040999999995, 0.60406223000000003, 0.57022952000000005, 0.57456998999999997, 0.61981388000000004, 0.64309300000000003, 0.75403390999999997, 0.75266655000000005, 0.65407048000000001, 0.52980380000000005, 0.79984374999999996, 0.68108995000000006, 0.61069605999999999, 0.79956828000000002, 0.68525809000000004, 0.77770578000000001, 0.77537922999999997, 0.76605701000000004, 0.70625411999999999, 0.68403398000000004, 0.54979301999999997, 0.57080569999999997, 0.45926551999999998, 0.53817833000000004, 0.45542998000000001, 0.52625016000000002, 0.48286651000000003, 0.51417950999999995, 0.78811969999999998, 0.68301993000000005, 0.79555133, 0.69489650000000003, 0.63523816, 0.73165648999999999, 0.60330846999999999, 0.73612842000000001, 0.61764971000000002, 0.56887445999999997, 0.49888159999999998, 0.49320227999999999, 0.63520029, 0.60815931999999995, 0.63034171000000006, 0.92437647999999994, 0.61072804000000003, 0.75338037999999996, 0.75603717000000004, 0.58176901999999997, 0.73021526999999997)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}

# Approximate the narrower "best fit" width Excel computed for the new column
$ws.Columns("A").ColumnWidth = 9

# Selection moved onto the newly inserted column in the saved file
$null = $ws.Range("A2:A51").Select()
